# "Fix mistakes on project 01"
# Adds Day 7 / Day 8 / Day 9 entries to the 30-day log, widens column C
# to fit the new (longer) text, and leaves the selection on B13 (the
# first empty cell below the new data), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Day 7 -----------------------------------------------------
# Copy D9's formatting first so the new date cell reuses the existing
# built-in date style (numFmtId 14) instead of Excel minting a fresh
# (duplicate) custom number format when a date value is assigned.
$ws.Range("D9").Copy($ws.Range("D10"))
$ws.Range("B10").Value = "Day 7"
$ws.Range("C10").Value = "Finish project 01 (Scrumble Computer Science Problems) and start arrays and linked lists."
$ws.Range("D10").Value = "9/24/2020"

# --- Row 11: Day 8 (no description / date yet) -------------------------
$ws.Range("B11").Value = "Day 8"

# --- Row 12: Day 9 (no description / date yet) -------------------------
$ws.Range("B12").Value = "Day 9"

# --- Widen column C so the longer "Day 7" text fits --------------------
$ws.Columns("C").ColumnWidth = 75.17

# --- Move the active selection to B13 -----------------------------------
$ws.Range("B13").Select() | Out-Null
